# Update planted column load to 250kN -> regenerated beam design results
# Values below reflect the recalculated Moment (E), Shear (F), Calc Depth (D)
# and Rebar (G) columns on the "Beams Analysis" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Beams Analysis")

# Row 2 (TR-1)
$ws.Range("D2").Value = 0.5
$ws.Range("E2").Value = 290
$ws.Range("F2").Value = 212
$ws.Range("G2").Value = "10 T16"

# Row 4 (TR-3)
$ws.Range("E4").Value = 81
$ws.Range("F4").Value = 185
$ws.Range("G4").Value = "4 T16"

# Row 7 (TR-6)
$ws.Range("D7").Value = 0.45
$ws.Range("E7").Value = 252
$ws.Range("F7").Value = 207
$ws.Range("G7").Value = "10 T16"

# Row 8 (TR-7)
$ws.Range("D8").Value = 0.4
$ws.Range("E8").Value = 167
$ws.Range("F8").Value = 196
$ws.Range("G8").Value = "7 T16"

# Row 9 (TR-8)
$ws.Range("E9").Value = 70
$ws.Range("F9").Value = 184
$ws.Range("G9").Value = "3 T16"

# Row 13 (TR-12)
$ws.Range("D13").Value = 0.55
$ws.Range("E13").Value = 401
$ws.Range("F13").Value = 226
$ws.Range("G13").Value = "12 T16"

# Row 22 (TR-21)
$ws.Range("E22").Value = 104
$ws.Range("F22").Value = 188
$ws.Range("G22").Value = "5 T16"

# Row 23 (TR-22)
$ws.Range("D23").Value = 0.4
$ws.Range("E23").Value = 196
$ws.Range("F23").Value = 200
$ws.Range("G23").Value = "9 T16"

# Row 30 (TR-29)
$ws.Range("E30").Value = 104
$ws.Range("F30").Value = 188
$ws.Range("G30").Value = "5 T16"

# Row 31 (TR-30)
$ws.Range("D31").Value = 0.55
$ws.Range("E31").Value = 401
$ws.Range("F31").Value = 226
$ws.Range("G31").Value = "12 T16"

# Row 32 (TR-31)
$ws.Range("D32").Value = 0.5
$ws.Range("E32").Value = 290
$ws.Range("F32").Value = 212
$ws.Range("G32").Value = "10 T16"

# Row 47 (TR-46)
$ws.Range("D47").Value = 0.5
$ws.Range("E47").Value = 295
$ws.Range("F47").Value = 213
$ws.Range("G47").Value = "10 T16"

# Row 48 (TR-47)
$ws.Range("D48").Value = 0.4
$ws.Range("E48").Value = 212
$ws.Range("F48").Value = 202
$ws.Range("G48").Value = "9 T16"

# Row 49 (TR-48)
$ws.Range("D49").Value = 0.4
$ws.Range("E49").Value = 167
$ws.Range("F49").Value = 196
$ws.Range("G49").Value = "8 T16"

# Row 50 (TR-49)
$ws.Range("E50").Value = 71
$ws.Range("F50").Value = 184
$ws.Range("G50").Value = "3 T16"

# Row 51 (TR-50)
$ws.Range("D51").Value = 0.45
$ws.Range("E51").Value = 269
$ws.Range("F51").Value = 209
$ws.Range("G51").Value = "10 T16"

# Row 52 (TR-51)
$ws.Range("D52").Value = 0.5
$ws.Range("E52").Value = 344
$ws.Range("F52").Value = 219
$ws.Range("G52").Value = "12 T16"

# Row 54 (TR-53)
$ws.Range("D54").Value = 0.5
$ws.Range("E54").Value = 344
$ws.Range("F54").Value = 219
$ws.Range("G54").Value = "12 T16"

# Row 55 (TR-54)
$ws.Range("D55").Value = 0.45
$ws.Range("E55").Value = 269
$ws.Range("F55").Value = 209
$ws.Range("G55").Value = "10 T16"

# Row 56 (TR-55)
$ws.Range("D56").Value = 0.4
$ws.Range("E56").Value = 212
$ws.Range("F56").Value = 202
$ws.Range("G56").Value = "9 T16"

# Row 57 (TR-56)
$ws.Range("D57").Value = 0.4
$ws.Range("E57").Value = 207
$ws.Range("F57").Value = 201
$ws.Range("G57").Value = "9 T16"

# Row 58 (TR-57)
$ws.Range("D58").Value = 0.4
$ws.Range("E58").Value = 207
$ws.Range("F58").Value = 201
$ws.Range("G58").Value = "9 T16"
